$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 41: a game session record appended to the "Electromecanica" log
# (5500 tiradas de juego electromecanica).

# --- Column A: date/time string (kept as literal text, not a date serial) ---
$ws.Range("A41").Value = "2024-10-14 22:02:12"

# --- Column B: empty text placeholder (matches the other empty inlineStr cells) ---
$ws.Range("B41").Value = "'"
$ws.Range("B41").Style = "Normal"

# --- Numeric columns ---
$ws.Range("C41").Value = 10
$ws.Range("D41").Value = 3
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 0

# --- Columns I-M: empty text placeholders ---
$ws.Range("I41").Value = "'"
$ws.Range("I41").Style = "Normal"
$ws.Range("J41").Value = "'"
$ws.Range("J41").Style = "Normal"
$ws.Range("K41").Value = "'"
$ws.Range("K41").Style = "Normal"
$ws.Range("L41").Value = "'"
$ws.Range("L41").Style = "Normal"
$ws.Range("M41").Value = "'"
$ws.Range("M41").Style = "Normal"

# --- More numeric columns ---
$ws.Range("N41").Value = 10
$ws.Range("O41").Value = 10
$ws.Range("P41").Value = 1

# --- Column Q: empty text placeholder ---
$ws.Range("Q41").Value = "'"
$ws.Range("Q41").Style = "Normal"

$ws.Range("R41").Value = 5

# --- Column S: empty text placeholder ---
$ws.Range("S41").Value = "'"
$ws.Range("S41").Style = "Normal"

$ws.Range("T41").Value = 20

# --- Column U: literal "25%" text (not an auto-converted percentage number) ---
$ws.Range("U41").NumberFormat = "@"
$ws.Range("U41").Value = "25%"
$ws.Range("U41").Style = "Normal"

# --- Column V: ruleta data file path ---
$ws.Range("V41").Value = "D:\Repositorio\jonatha1992\Predictor_ruleta\Data\Electromecanica.xlsx"

# --- Column W: empty text placeholder ---
$ws.Range("W41").Value = "'"
$ws.Range("W41").Style = "Normal"

# --- Column X: simulation flag text ---
$ws.Range("X41").Value = "No es Simulación"

# --- Column Y: predicted numbers count ---
$ws.Range("Y41").Value = 40
